# Update for DWSIM new json format
#
# "Pump Info" sheet (Worksheets.Item(1)) gains two new rows describing the
# pump curve ("Name" / "Description"), and the "Pump Data" sheet
# (Worksheets.Item(2)) switches its unit-of-measure header row from
# m3/s, m, kW, abs -> m3/h, ft, HP, %.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Pump Info"
$ws2 = $wb.Worksheets.Item(2)   # "Pump Data"

# --- Pump Data: change the unit-of-measure row (row 2) -------------------
$ws2.Range("A2").Value = "m3/h"
$ws2.Range("B2").Value = "ft"
$ws2.Range("C2").Value = "HP"
$ws2.Range("D2").Value = "%"
$ws2.Range("E2").Value = "m"

# --- Pump Info: add Name / Description rows -------------------------------
# Fill column-by-column (A3, A4, then B3, B4) so that new shared strings are
# created in the same order as the authored workbook.
$ws1.Range("A3").Value = "Name"
$ws1.Range("A4").Value = "Description"
$ws1.Range("B3").Value = "Sulzer Pump Curve"
$ws1.Range("B4").Value = "Sulzer Pump Curve Test"

# Widen column B on the Pump Info sheet to fit the new values.
$ws1.Columns.Item(2).ColumnWidth = 11.333333

# --- View state: make "Pump Info" the active/selected tab ----------------
# First leave the last selection on "Pump Data" where it was left (F6),
# then switch to and activate "Pump Info".
$ws2.Activate()
$ws2.Range("F6").Select()
$ws1.Activate()
